$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Fold_1")
$ws.Range("B2").Value = 6.87565
$ws.Range("C2").Value = 1.2232
$ws.Range("D2").Value = 35.04514999999999
$ws.Range("E2").Value = 0.0359
$ws.Range("F2").Value = 25.56695
$ws.Range("G2").Value = 25.56695
$ws.Range("K2").Value = 65.9374
$ws.Range("L2").Value = 25.567
$ws.Range("M2").Value = 40.3704
$ws.Range("N2").Value = 7.785999999999999
$ws.Range("O2").Value = 32.5844
$ws.Range("B3").Value = 16.949
$ws.Range("C3").Value = 3.07
$ws.Range("D3").Value = 53.744
$ws.Range("F3").Value = 28.401
$ws.Range("G3").Value = 27.349
$ws.Range("K3").Value = 28.334
$ws.Range("L3").Value = 27.349
$ws.Range("M3").Value = 0.985
$ws.Range("N3").Value = 0.985
$ws.Range("B4").Value = 33.529
$ws.Range("C4").Value = 2.449
$ws.Range("D4").Value = 53.46
$ws.Range("F4").Value = 35.893
$ws.Range("G4").Value = 35.8925808219178
$ws.Range("K4").Value = 29.259
$ws.Range("L4").Value = 29.259

$ws = $wb.Worksheets.Item("Fold_2")
$ws.Range("B2").Value = 6.19875
$ws.Range("C2").Value = 1.20455
$ws.Range("D2").Value = 34.63824999999999
$ws.Range("E2").Value = 0.0359
$ws.Range("F2").Value = 25.47575
$ws.Range("G2").Value = 25.47575
$ws.Range("K2").Value = 228.6622
$ws.Range("L2").Value = 25.476
$ws.Range("M2").Value = 203.1862
$ws.Range("N2").Value = 8.9472
$ws.Range("O2").Value = 194.2392
$ws.Range("B3").Value = 14.045
$ws.Range("C3").Value = 3.349
$ws.Range("D3").Value = 53.613
$ws.Range("F3").Value = 27.519
$ws.Range("G3").Value = 27.029
$ws.Range("H3").Value = 0.489
$ws.Range("I3").Value = 0.489
$ws.Range("K3").Value = 32.8484
$ws.Range("L3").Value = 27.029
$ws.Range("M3").Value = 5.8194
$ws.Range("N3").Value = 5.8194
$ws.Range("B4").Value = 19.899
$ws.Range("C4").Value = 2.786
$ws.Range("D4").Value = 53.876
$ws.Range("F4").Value = 31.612
$ws.Range("G4").Value = 31.61164383561644
$ws.Range("K4").Value = 31.134
$ws.Range("L4").Value = 27.675
$ws.Range("M4").Value = 3.459000000000001
$ws.Range("N4").Value = 3.459000000000001

$ws = $wb.Worksheets.Item("Fold_3")
$ws.Range("B2").Value = 7.723350000000001
$ws.Range("C2").Value = 1.3063
$ws.Range("D2").Value = 32.7205
$ws.Range("F2").Value = 25.62845
$ws.Range("G2").Value = 25.62845
$ws.Range("K2").Value = 75.7448
$ws.Range("L2").Value = 25.628
$ws.Range("M2").Value = 50.1168
$ws.Range("N2").Value = 5.8988
$ws.Range("O2").Value = 44.218
$ws.Range("B3").Value = 19.899
$ws.Range("C3").Value = 2.786
$ws.Range("D3").Value = 53.876
$ws.Range("F3").Value = 28.539
$ws.Range("G3").Value = 27.675
$ws.Range("H3").Value = 0.865
$ws.Range("I3").Value = 0.865
$ws.Range("K3").Value = 27.675
$ws.Range("L3").Value = 27.675
$ws.Range("B4").Value = 33.529
$ws.Range("C4").Value = 2.449
$ws.Range("D4").Value = 53.46
$ws.Range("F4").Value = 35.893
$ws.Range("G4").Value = 35.8925808219178
$ws.Range("K4").Value = 29.259
$ws.Range("L4").Value = 29.259

$ws = $wb.Worksheets.Item("Fold_4")
$ws.Range("B2").Value = 8.298950000000001
$ws.Range("C2").Value = 1.2959
$ws.Range("D2").Value = 32.83895
$ws.Range("E2").Value = 0.0359
$ws.Range("F2").Value = 25.6993
$ws.Range("G2").Value = 25.6993
$ws.Range("K2").Value = 45.387
$ws.Range("L2").Value = 25.699
$ws.Range("M2").Value = 19.688
$ws.Range("N2").Value = 2.7378
$ws.Range("O2").Value = 16.9502
$ws.Range("B3").Value = 19.899
$ws.Range("C3").Value = 2.786
$ws.Range("D3").Value = 53.876
$ws.Range("F3").Value = 28.539
$ws.Range("G3").Value = 27.675
$ws.Range("H3").Value = 0.865
$ws.Range("I3").Value = 0.865
$ws.Range("K3").Value = 27.675
$ws.Range("L3").Value = 27.675
$ws.Range("B4").Value = 33.529
$ws.Range("C4").Value = 2.449
$ws.Range("D4").Value = 53.46
$ws.Range("F4").Value = 35.893
$ws.Range("G4").Value = 35.8925808219178
$ws.Range("K4").Value = 29.259
$ws.Range("L4").Value = 29.259

$ws = $wb.Worksheets.Item("Fold_5")
$ws.Range("B2").Value = 8.0367
$ws.Range("C2").Value = 1.05765
$ws.Range("D2").Value = 33.68215
$ws.Range("E2").Value = 0.0359
$ws.Range("F2").Value = 25.66015000000001
$ws.Range("G2").Value = 25.66015000000001
$ws.Range("K2").Value = 144.3206
$ws.Range("L2").Value = 25.66
$ws.Range("M2").Value = 118.6606
$ws.Range("N2").Value = 3.2118
$ws.Range("O2").Value = 115.4488
$ws.Range("B3").Value = 19.917
$ws.Range("C3").Value = 2.49
$ws.Range("D3").Value = 52.862
$ws.Range("F3").Value = 28.489
$ws.Range("G3").Value = 27.625
$ws.Range("H3").Value = 0.865
$ws.Range("I3").Value = 0.865
$ws.Range("K3").Value = 30.8814
$ws.Range("L3").Value = 27.625
$ws.Range("M3").Value = 3.2564
$ws.Range("N3").Value = 0.1934
$ws.Range("O3").Value = 3.0628
$ws.Range("B4").Value = 33.529
$ws.Range("C4").Value = 2.449
$ws.Range("D4").Value = 52.947
$ws.Range("F4").Value = 35.882
$ws.Range("G4").Value = 35.88203287671233
$ws.Range("K4").Value = 29.292
$ws.Range("L4").Value = 29.248
$ws.Range("M4").Value = 0.044
$ws.Range("N4").Value = 0.044
$ws.Range("O4").Value = 0
